# AutoCommit_21 мая 2024 г. 13:48:16_SibNout2023
# Bump F/G/H scores to 5 for two students (rows 5 and 12), which also
# updates their SUM() totals in column J, and moves the frozen-pane
# view/selection up to row 5 (from row 15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("Ефимов Владислав"): F/G/H -> 5/5/5 (was 2/blank/blank)
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 5

# Row 12 ("Кравчук Мария"): F/G/H -> 5/5/5 (was 2/blank/blank)
$ws.Range("F12").Value = 5
$ws.Range("G12").Value = 5
$ws.Range("H12").Value = 5

# Scroll/select so the active cell in the frozen bottom-right pane is I5
# (matches the updated <selection pane="bottomRight" activeCell="I5" .../>)
$ws.Range("I5").Select()
